# Remove the "syscall" / "break" rows from the R-Type (Funct) table.
# These were rows 10 and 11 on the "R-Type" worksheet; deleting the
# entire rows shifts the remaining Funct-table rows up and shrinks
# both the worksheet dimension and Table2's range automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R-Type")

$ws.Range("A10:D11").EntireRow.Delete()

# Update the view's active selection to reflect where the user ended up
# after the edit (cell E12, with no frozen/top-left scroll override).
$ws.Range("E12").Select()
